$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Move the footer block (rows 34-35) down to rows 42-43 so that
#    there is room to grow the data table (which now needs rows
#    16-37 instead of 16-29).
# ------------------------------------------------------------------
$ws.Range("B34:C34").UnMerge()
$ws.Range("B35:C35").UnMerge()
$ws.Range("H34:J34").UnMerge()
$ws.Range("H35:J35").UnMerge()

$ws.Range("B34:J35").Cut($ws.Range("B42:J43"))
$ws.Rows("34:35").Clear()
$ws.Range("D42:G43").Clear()

$ws.Range("B42:C42").Merge()
$ws.Range("B43:C43").Merge()
$ws.Range("H42:J42").Merge()
$ws.Range("H43:J43").Merge()

# ------------------------------------------------------------------
# 2) Grow the data table. Row 29 currently carries the "last row"
#    (bottom border) style; duplicate that style onto the new last
#    row (37) first, then restyle rows 29-36 as normal middle rows
#    by copying the formatting of row 28.
# ------------------------------------------------------------------
$ws.Range("B29:J29").Copy($ws.Range("B37:J37"))

for ($r = 29; $r -le 36; $r++) {
  $ws.Range("B28:J28").Copy($ws.Range("B" + $r + ":J" + $r))
}

# ------------------------------------------------------------------
# 3) Populate the table values.
#    IVAN DAVID GOMEZ MENDOZA now covers periods 2406-2507 (14 rows)
#    CHRISTIAN FLOREZ CANABAL now covers periods 2407-2502 (8 rows)
# ------------------------------------------------------------------
$ivan = @("2507","2506","2505","2504","2503","2502","2501","2412","2411","2410","2409","2408","2407","2406")
$row = 16
foreach ($p in $ivan) {
  $ws.Range("B" + $row).Value = "CC"
  $ws.Range("C" + $row).Value = "8852831"
  $ws.Range("D" + $row).Value = "IVAN DAVID GOMEZ MENDOZA"
  $ws.Range("E" + $row).Value = $p
  $ws.Range("F" + $row).Value = 75467
  $ws.Range("G" + $row).Value = 1886679
  $row = $row + 1
}

$christian = @("2502","2501","2412","2411","2410","2409","2408","2407")
foreach ($p in $christian) {
  $ws.Range("B" + $row).Value = "CC"
  $ws.Range("C" + $row).Value = "73184636"
  $ws.Range("D" + $row).Value = "CHRISTIAN FLOREZ CANABAL"
  $ws.Range("E" + $row).Value = $p
  $ws.Range("F" + $row).Value = 63554
  $ws.Range("G" + $row).Value = 1588857
  $row = $row + 1
}

# ------------------------------------------------------------------
# 4) Update the summary figures above the table.
#    Valor Mora total = 14*75467 + 8*63554 = 1564970
#    Cant. Periodos = 14
# ------------------------------------------------------------------
$ws.Range("E11").Value = 1564970
$ws.Range("F13").Value = 14

Write-Host "done"
